$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the cached "datetimeFigureOut" date field text that lives on
#    the slide master and every slide layout (Insert > Header & Footer ->
#    Date automatically updates; PowerPoint re-caches the literal text
#    shown for the field whenever the deck is saved).
# ---------------------------------------------------------------------
function Set-DatePlaceholderText($container) {
    $shapeCount = $container.Shapes.Count
    for ($i = 1; $i -le $shapeCount; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.Type -eq 14) {
            $isDatePlaceholder = $false
            try {
                if ($shp.PlaceholderFormat.Type -eq 16) {
                    $isDatePlaceholder = $true
                }
            } catch {
            }
            if ($isDatePlaceholder -and $shp.HasTextFrame) {
                if ($shp.TextFrame.TextRange.Text -eq "10/2/2020") {
                    $shp.TextFrame.TextRange.Text = "10/22/2020"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Set-DatePlaceholderText($master)

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Set-DatePlaceholderText($layouts.Item($li))
}

# ---------------------------------------------------------------------
# 2) Fix the typo on Slide 2 ("Wheel Encoders - convert tick to rpm"):
#    "...40 ticks per resolution" -> "...40 ticks per revolution"
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$shapeCount2 = $slide2.Shapes.Count
for ($i = 1; $i -le $shapeCount2; $i++) {
    $shp = $slide2.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        if ($tr.Text -eq "Our wheel encoder has 40 ticks per resolution") {
            $sub = $tr.Characters(32, 14)
            $sub.Text = "per revolution"
        }
    }
}
